$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.624.68'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").Value = '1.898.82'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.66%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.43'
$ws.Range("E5").Value = '  -2.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.693'
$ws.Range("E6").Value = '  -4.36%  '
$ws.Range("E7").Value = '  -0.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.04'
$ws.Range("E8").Value = '  +8.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.353'
$ws.Range("E9").Value = '  -6.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0742'
$ws.Range("E10").Value = '  -2.43%  '
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '13.12'
$ws.Range("E12").Value = '  +1.31%  '
$ws.Range("D13").Value = '2.177.22'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.731'
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.96'
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("D16").Value = '1.898.91'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '35.649.37'
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '73.88'
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("D19").Value = '0.0₃0826'
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '247.58'
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.92'
$ws.Range("E21").Value = '  -0.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.98'
$ws.Range("E22").Value = '  -2.54%  '
$ws.Range("E23").Value = '  -0.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.55'
$ws.Range("E24").Value = '  +3.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  -9.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.88'
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.44'
$ws.Range("E28").Value = '  -1.75%  '
$ws.Range("E29").Value = '  -3.82%  '
$ws.Range("E31").Value = '  +7.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.27'
$ws.Range("E32").Value = '  -2.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0583'
$ws.Range("E33").Value = '  -0.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.24'
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.855'
$ws.Range("E36").Value = '  -6.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.02'
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.55'
$ws.Range("E38").Value = '  -22.26%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0684'
$ws.Range("E39").Value = '  +5.22%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.20'
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.10'
$ws.Range("E41").Value = '  +1.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0214'
$ws.Range("E42").Value = '  -2.08%  '
$ws.Range("E43").Value = '  -1.99%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.301.12'
$ws.Range("E44").Value = '  -2.74%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.38'
$ws.Range("E45").Value = '  -2.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0807'
$ws.Range("E46").Value = '  +6.81%  '
$ws.Range("E47").Value = '  -1.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.09'
$ws.Range("E49").Value = '  +3.39%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.38'
$ws.Range("E50").Value = '  -4.84%  '
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.45'
$ws.Range("E51").Value = '  -3.68%  '
